$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Vase"
$ws.Range("B2").Value = "Cat"
$ws.Range("C2").Value = "Ukraine"
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 120
$ws.Range("F2").Value = 1080

$ws.Range("F3").Value = 1080
